$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("O2").Value = 0.06816352613805679
$ws.Range("P2").Value = 0.0681635261380568
$ws.Range("Q2").Value = 0.01327427880666666
$ws.Range("S2").Value = 0.06816352613805679
$ws.Range("T2").Value = 0.0681635261380568

# Row 3 updates
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("M3").Value = 1.873282666666666
$ws.Range("N3").Value = 5.619847999999999
$ws.Range("O3").Value = 0.9318364738619431
$ws.Range("P3").Value = 0.9318364738619432
$ws.Range("Q3").Value = 0.1814673896302222
$ws.Range("R3").Value = 1.633206506672
$ws.Range("S3").Value = 0.9318364738619431
$ws.Range("T3").Value = 0.9318364738619432

# Remove the now-obsolete rows 4 and 5 (their data moved into / was
# superseded by the updated rows 2-3 above).
$ws.Rows("4:5").Delete()
